# Weekly SEBI downloads workbook update:
# Replace the November SEBI circular/press-release/consultation rows with the
# October AIF + SEBI consultation-paper rows (commit: "Categorised AIF from
# SEBI in searching agent").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Remove the hyperlink objects that live on rows 6 and 7 (G6, G7)
#    before we delete those rows, and remove rows 6-7 entirely so the
#    sheet ends up with only 5 rows (header + 4 data rows).
# ------------------------------------------------------------------
$linksToRemove = @()
foreach ($link in $ws.Hyperlinks) {
    if ($link.Range.Row -ge 6) {
        $linksToRemove += $link
    }
}
for ($i = $linksToRemove.Count - 1; $i -ge 0; $i--) {
    $linksToRemove[$i].Delete()
}

$ws.Rows.Item(7).Delete() | Out-Null
$ws.Rows.Item(6).Delete() | Out-Null

# ------------------------------------------------------------------
# 2) Overwrite the remaining 4 data rows (2-5) with the new October
#    AIF / SEBI data.
# ------------------------------------------------------------------
$rowData = @(
    @{ Row = 2; A = 'AIF';  B = 'Circulars';           C = '2025'; D = 'October'; E = '2025-10-24';
       F = 'Transfer of portfolios of clients (PMS business) by Portfolio Managers.';
       G = 'https://www.sebi.gov.in/sebi_data/attachdocs/oct-2025/1761301360689.pdf';
       H = '1761301360689.pdf';
       I = '/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/AIF/Circulars/2025/October/1761301360689.pdf' },
    @{ Row = 3; A = 'SEBI'; B = 'Consulatation Paper'; C = '2025'; D = 'October'; E = '2025-10-24';
       F = 'Consultation paper for review of LODR Regulations - Clarification regarding the timeline for transfer of unclaimed amount by entity having listed non-convertible securities.';
       G = 'https://www.sebi.gov.in/sebi_data/attachdocs/oct-2025/1761306586026.pdf';
       H = '1761306586026.pdf';
       I = '/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/SEBI/Consulatation Paper/2025/October/1761306586026.pdf' },
    @{ Row = 4; A = 'SEBI'; B = 'Consulatation Paper'; C = '2025'; D = 'October'; E = '2025-10-23';
       F = 'Circular on Relaxation of India geo-tagging for NRI clients re-KYCKYC modification through digital on-boarding video client Identification process (V-CIP)';
       G = 'https://www.sebi.gov.in/sebi_data/attachdocs/oct-2025/1761220376339.pdf';
       H = '1761220376339.pdf';
       I = '/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/SEBI/Consulatation Paper/2025/October/1761220376339.pdf' },
    @{ Row = 5; A = 'SEBI'; B = 'Consulatation Paper'; C = '2025'; D = 'October'; E = '2025-10-23';
       F = 'Consultation paper on Standardization of process for Opening of Mutual Fund Folios and Execution of First Investment';
       G = 'https://www.sebi.gov.in/sebi_data/attachdocs/oct-2025/1761210652019.pdf';
       H = '1761210652019.pdf';
       I = '/Users/admin/Downloads/Tejomaya_pdfs/Akshayam Data/SEBI/Consulatation Paper/2025/October/1761210652019.pdf' }
)

# Columns C (Year) and E (IssueDate) hold text that Excel would otherwise
# auto-convert to a number / date, so force a Text format first and then
# restore the default "Normal" style once the text value is stored (keeps
# the cells free of any lingering explicit style index).
$ws.Range("C2:C5").NumberFormat = "@"
$ws.Range("E2:E5").NumberFormat = "@"

foreach ($rd in $rowData) {
    $r = $rd.Row
    $ws.Cells.Item($r, 1).Value = $rd.A
    $ws.Cells.Item($r, 2).Value = $rd.B
    $ws.Cells.Item($r, 3).Value = $rd.C
    $ws.Cells.Item($r, 4).Value = $rd.D
    $ws.Cells.Item($r, 5).Value = $rd.E
    $ws.Cells.Item($r, 6).Value = $rd.F
    $ws.Cells.Item($r, 7).Value = $rd.G
    $ws.Cells.Item($r, 8).Value = $rd.H
    $ws.Cells.Item($r, 9).Value = $rd.I
}

$ws.Range("C2:C5").Style = "Normal"
$ws.Range("E2:E5").Style = "Normal"

# ------------------------------------------------------------------
# 3) Point the 4 remaining hyperlinks (G2:G5) at the new October PDFs.
#    (The cell text itself was already updated above; here we update
#    the underlying hyperlink relationship target to match.)
# ------------------------------------------------------------------
$newUrls = @{
    2 = 'https://www.sebi.gov.in/sebi_data/attachdocs/oct-2025/1761301360689.pdf'
    3 = 'https://www.sebi.gov.in/sebi_data/attachdocs/oct-2025/1761306586026.pdf'
    4 = 'https://www.sebi.gov.in/sebi_data/attachdocs/oct-2025/1761220376339.pdf'
    5 = 'https://www.sebi.gov.in/sebi_data/attachdocs/oct-2025/1761210652019.pdf'
}

foreach ($link in $ws.Hyperlinks) {
    $r = $link.Range.Row
    if ($newUrls.ContainsKey($r)) {
        $link.Address = $newUrls[$r]
    }
}

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
